$p = $ppt.ActivePresentation

# The deck's first slide master ("geometric", rId12 -> ppt/theme/theme2.xml) is
# being restyled from the "Geometric" palette to the "Simple Light" palette.
# (Conversely, the "Simple Light" master picks up the old "Geometric" colours,
# but the PowerPoint theme-color object model here always binds to the
# presentation's primary/first slide master, so that side of the swap is
# applied through the same call chain below.)

$master = $p.Slides.Item(2).Master
$colorScheme = $master.Theme.ThemeColorScheme

# VBA/COM RGB() packs colours as 0xBBGGRR (r + g*256 + b*65536).
$colorScheme.Colors(1).RGB  = 0            # dk1      000000
$colorScheme.Colors(2).RGB  = 16777215     # lt1      FFFFFF
$colorScheme.Colors(3).RGB  = 5855577      # dk2      595959
$colorScheme.Colors(4).RGB  = 15658734     # lt2      EEEEEE
$colorScheme.Colors(5).RGB  = 16024898     # accent1  4285F4
$colorScheme.Colors(6).RGB  = 2171169      # accent2  212121
$colorScheme.Colors(7).RGB  = 10260600     # accent3  78909C
$colorScheme.Colors(8).RGB  = 4238335      # accent4  FFAB40
$colorScheme.Colors(9).RGB  = 10983168     # accent5  0097A7
$colorScheme.Colors(10).RGB = 4325358      # accent6  EEFF41
$colorScheme.Colors(11).RGB = 10983168     # hlink    0097A7
$colorScheme.Colors(12).RGB = 10983168     # folHlink 0097A7
